$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '29.402.23'
$ws.Range("E2").NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '  -0.16%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.850.18'
$ws.Range("E3").NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '  -0.02%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.0000'
$ws.Range("E4").NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '240.60'
$ws.Range("E5").NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '  -0.15%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.6290'
$ws.Range("E6").NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '  -0.21%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.001'
$ws.Range("E7").NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '  +0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.07593'
$ws.Range("E8").NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '  -1.37%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.2922'
$ws.Range("E9").NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '  -0.27%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '24.54'
$ws.Range("E10").NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '  -0.72%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07755'
$ws.Range("E11").NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '  +0.17%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.849.56'
$ws.Range("E12").NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '  -2.10%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.010'
$ws.Range("E13").NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '  -0.52%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.6793'
$ws.Range("E14").NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = '  -0.06%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.00001041'
$ws.Range("E15").NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '  -3.00%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '83.15'
$ws.Range("E16").NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '  -0.67%  '

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Cells.Item(17, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").NumberFormat = "@"
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '2.098.24'
$ws.Range("E17").NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '  -2.66%  '

# Row 18
$ws.Range("B18").NumberFormat = "@"
$ws.Cells.Item(18, 2).Value = 'Uniswap'
$ws.Range("C18").NumberFormat = "@"
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '6.108'
$ws.Range("E18").NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '  -1.55%  '

# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Cells.Item(19, 2).Value = 'WrappedBTC'
$ws.Range("C19").NumberFormat = "@"
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '29.390.10'
$ws.Range("E19").NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '  -0.26%  '

# Row 20
$ws.Range("B20").NumberFormat = "@"
$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Range("C20").NumberFormat = "@"
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '229.74'
$ws.Range("E20").NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '  +0.53%  '

# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Cells.Item(21, 2).Value = 'Avalanche'
$ws.Range("C21").NumberFormat = "@"
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '12.34'
$ws.Range("E21").NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = '  -0.89%  '

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Cells.Item(22, 2).Value = 'Dai'
$ws.Range("C22").NumberFormat = "@"
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '1.001'
$ws.Range("E22").NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '  +0.06%  '

# Row 23
$ws.Range("B23").NumberFormat = "@"
$ws.Cells.Item(23, 2).Value = 'Chainlink'
$ws.Range("C23").NumberFormat = "@"
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '7.435'
$ws.Range("E23").NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '  -0.34%  '

# Row 24
$ws.Range("B24").NumberFormat = "@"
$ws.Cells.Item(24, 2).Value = 'BinanceUSD'
$ws.Range("C24").NumberFormat = "@"
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '1.000'
$ws.Range("E24").NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '  -0.02%  '

# Row 25
$ws.Range("B25").NumberFormat = "@"
$ws.Cells.Item(25, 2).Value = 'Monero'
$ws.Range("C25").NumberFormat = "@"
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '159.28'
$ws.Range("E25").NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '  +1.12%  '

# Row 26
$ws.Range("B26").NumberFormat = "@"
$ws.Cells.Item(26, 2).Value = 'Stellar'
$ws.Range("C26").NumberFormat = "@"
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.1396'
$ws.Range("E26").NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '  +1.00%  '

# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Cells.Item(27, 2).Value = 'Cosmos'
$ws.Range("C27").NumberFormat = "@"
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '8.452'
$ws.Range("E27").NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '  +0.47%  '

# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Range("C28").NumberFormat = "@"
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '17.67'
$ws.Range("E28").NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = '  -0.17%  '

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Range("C29").NumberFormat = "@"
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.425'
$ws.Range("E29").NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = '  +5.79%  '

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Range("C30").NumberFormat = "@"
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.474'
$ws.Range("E30").NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = '  +0.28%  '

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Cells.Item(31, 2).Value = 'Hedera'
$ws.Range("C31").NumberFormat = "@"
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.05690'
$ws.Range("E31").NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = '  +0.20%  '

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Cells.Item(32, 2).Value = 'Filecoin'
$ws.Range("C32").NumberFormat = "@"
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.119'
$ws.Range("E32").NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = '  -0.29%  '

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Cells.Item(33, 2).Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").NumberFormat = "@"
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.045'
$ws.Range("E33").NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = '  +0.19%  '

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Cells.Item(34, 2).Value = 'ARBITRUM'
$ws.Range("C34").NumberFormat = "@"
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.157'
$ws.Range("E34").NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = '  -0.64%  '

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Cells.Item(35, 2).Value = 'LidoDAOToken'
$ws.Range("C35").NumberFormat = "@"
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.824'
$ws.Range("E35").NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = '  -1.47%  '

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Cells.Item(36, 2).Value = 'ImmutableX'
$ws.Range("C36").NumberFormat = "@"
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.6964'
$ws.Range("E36").NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = '  -1.06%  '

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Cells.Item(37, 2).Value = 'HuobiToken'
$ws.Range("C37").NumberFormat = "@"
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.583'
$ws.Range("E37").NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = '  -0.12%  '

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Range("C38").NumberFormat = "@"
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01826'
$ws.Range("E38").NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = '  +1.92%  '

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Cells.Item(39, 2).Value = 'Maker'
$ws.Range("C39").NumberFormat = "@"
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.240.04'
$ws.Range("E39").NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '  +1.66%  '

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Cells.Item(40, 2).Value = 'MXToken'
$ws.Range("C40").NumberFormat = "@"
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.718'
$ws.Range("E40").NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '  -2.27%  '

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Cells.Item(41, 2).Value = 'FraxShare'
$ws.Range("C41").NumberFormat = "@"
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '6.415'
$ws.Range("E41").NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '  -2.16%  '

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Range("C42").NumberFormat = "@"
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.9021'
$ws.Range("E42").NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '  -0.50%  '

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Cells.Item(43, 2).Value = 'PaxDollar'
$ws.Range("C43").NumberFormat = "@"
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.000'
$ws.Range("E43").NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '  -0.07%  '

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Cells.Item(44, 2).Value = 'RocketPoolETH'
$ws.Range("C44").NumberFormat = "@"
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.005.40'
$ws.Range("E44").NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '  -2.76%  '

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Cells.Item(45, 2).Value = 'Quant'
$ws.Range("C45").NumberFormat = "@"
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '101.47'
$ws.Range("E45").NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '  -0.26%  '

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Cells.Item(46, 2).Value = 'Aave'
$ws.Range("C46").NumberFormat = "@"
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '65.75'
$ws.Range("E46").NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '  -0.99%  '

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Cells.Item(47, 2).Value = 'Aptos'
$ws.Range("C47").NumberFormat = "@"
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '7.139'
$ws.Range("E47").NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '  +0.09%  '

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Cells.Item(48, 2).Value = 'TheSandbox'
$ws.Range("C48").NumberFormat = "@"
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D48").NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.3996'
$ws.Range("E48").NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '  -0.68%  '

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Cells.Item(49, 2).Value = 'BabyDogeCoin'
$ws.Range("C49").NumberFormat = "@"
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.00000000116'
$ws.Range("E49").NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '  -2.74%  '

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Cells.Item(50, 2).Value = 'Algorand'
$ws.Range("C50").NumberFormat = "@"
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.1156'
$ws.Range("E50").NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '  +1.10%  '

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Range("C51").NumberFormat = "@"
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '9.029'
$ws.Range("E51").NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = '  +0.15%  '
